$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 8: "Exchange that Shitty Printer" / "SATURDAY?"  ->  "Linked-in => add everybody from the Red Poole (add Alyssa Liddle)" / "IN PROGRESS"
# Copy the format used at row 7 (yellow "in progress" style) onto row 8 first.
[void]$ws.Range("A7:B7").Copy()
[void]$ws.Range("A8").PasteSpecial(-4122)
$ws.Range("A8").Value2 = "Linked-in => add everybody from the Red Poole (add Alyssa Liddle)"
$ws.Range("B8").Value2 = "IN PROGRESS"

# --- Row 9: "Linked-in => add everybody..." / "SATURDAY?"  ->  "Connect with Marc Johnson" / "TODO"  (style unchanged)
$ws.Range("A9").Value2 = "Connect with Marc Johnson"
$ws.Range("B9").Value2 = "TODO"

# --- Row 10: "Jeff, Keith, Phil..." / "WHENEVS" -> "Jeff Benson, Organize my Contacts Database => ..." / "TODO" (+ empty C10 cell)
# B10 takes on the style currently used by B11 ("TODO" style) before row 11 disappears.
[void]$ws.Range("B11").Copy()
[void]$ws.Range("B10").PasteSpecial(-4122)
$ws.Range("A10").Value2 = "Jeff Benson, Organize my Contacts Database =>  Keith, Phil, Andy, Brad, Andrei, Ajit, Ryan"
$ws.Range("B10").Value2 = "TODO"
$c10 = $ws.Range("C10")
$c10.HorizontalAlignment = -4152
$c10.Font.Bold = $false
$c10.Interior.Pattern = -4142

# --- Row 11 ("Organize my Contacts Database" / "TODO") is removed outright; everything below shifts up one row.
[void]$ws.Rows("11").Delete()

# After the shift: row 13 = Personal/Household (correct already), row 14 = Honda Accord/TODO (native s25/s20),
# row 15 = Download YNAB/TODO (native s16/s22). Target wants the "Download" content promoted to row 13 with
# the yellow "in progress" style, and Honda Accord kept at row 14 with its original (untouched) style.
# Insert a fresh blank row at 13 -- this pushes Honda Accord/Download back down by one WITHOUT touching their
# native formatting (Honda Accord lands back on row 14 exactly as it was).
[void]$ws.Rows("13").Insert()
[void]$ws.Range("A7:B7").Copy()
[void]$ws.Range("A13").PasteSpecial(-4122)
$ws.Range("A13").Value2 = "Download and start using YNAB"
$ws.Range("B13").Value2 = "IN PROGRESS"
# The original "Download" row (now pushed down to row 15) is now a duplicate -- remove it so everything below
# collapses back up to its target position.
[void]$ws.Rows("15").Delete()

# --- Column A width & active selection
$ws.Columns("A").ColumnWidth = 86.83333333333334
[void]$ws.Range("A7").Select()
